# uren.xlsx: Mr. Black is replaced by Mr. White in the first (B/C) column,
# and the freed-up second (E/F) column slot is now used by a new
# colleague, Mr. Pink. Mr. Blue (H/I) is unaffected. The tracking grid is
# also extended from row 10 down to row 38.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# --- Header row (row 2): who each column pair belongs to -------------------
$ws.Range("B2").Value = "Uren Mr. White"
$ws.Range("E2").Value = "Uren Mr. Pink"
# H2 ("Uren Mr. Blue") is unchanged.

# --- Summary table (column L, rows 5-7): row order mirrors B/E/H ----------
# L7 stays "Mr. Blue".
$ws.Range("L5").Value = "Mr. White"
$ws.Range("L6").Value = "Mr. Pink"

# --- Hour entries -----------------------------------------------------------
$ws.Range("C4").Value = 5
$ws.Range("F4").Value = 5

# New entry for Mr. Pink (column E/F): 28 May 2012, no hours filled in yet.
$ws.Range("B4").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("E5").Value = 41057

# --- Extend the blank, date-formatted grid down through row 38 -------------
$ws.Range("B4").Copy()
$ws.Range("B7:B38").PasteSpecial(-4122)
$ws.Range("E6:E36").PasteSpecial(-4122)
$ws.Range("H4:H37").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Selection --------------------------------------------------------------
$ws.Range("F5").Select()
